# #5: property boat&car done
#
# The "汽車" (car) sheet only had an ad-hoc 6-column layout (name, capacity,
# owner, register_date, register_reason, acquire_value). Bring it in line
# with the standard property-declaration layout used by the other sheets
# in this workbook by adding the property_category / category / date /
# legislator_name / legislator_id / source_file / index columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# The "date" column holds a literal yyyy-mm-dd text value ("2012-04-27")
# that Excel would otherwise silently reinterpret as a real date serial
# number. Force the cell to stay text before assigning it.
$ws.Cells.Item(2, 10).NumberFormat = "@"

# --- Row 1: header labels (B1:N1) ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2: the single car record (B2:N2) ---
$ws.Cells.Item(2, 2).Value = "自用小汽車（TOYOTA型號：ALTIS)"
$ws.Cells.Item(2, 3).Value = 1800
$ws.Cells.Item(2, 4).Value = "王琴賀"
$ws.Cells.Item(2, 5).Value = "97年02月29日"
$ws.Cells.Item(2, 6).Value = "買賣"
$ws.Cells.Item(2, 7).Value = 850000
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2012-04-27"
$ws.Cells.Item(2, 11).Value = "蔡煌瑯"
$ws.Cells.Item(2, 12).Value = 752
$ws.Cells.Item(2, 13).Value = "tmpd4981"
$ws.Cells.Item(2, 14).Value = 41

# --- Make the newly-added cells (H:N) pick up the same look as the rest
#     of the header / data row (bold+bordered header style, plain data
#     style) by copying the formatting across from the existing cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$wb.Save()
